$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column M (13) so the existing M column (health)
# shifts right to become column N, carrying its formulas/width along.
$ws.Columns.Item(13).Insert()

# Re-establish N3:N5 (the old M3:M5, shifted right by the insert) as a shared
# formula group - the column insert split it into standalone formulas.
$ws.Range("N3:N5").Formula = "=K3*0.001"

# Header for the newly inserted column M
$ws.Range("M1").Value = "short_percent_change_km_yr"

# Fill formulas for the new column M (L/1000), matching the existing row pattern:
# M2 stands alone, M3:M5 share one formula group (mirrors the original M column's shape).
$ws.Range("M2").Formula = "=L2/1000"
$ws.Range("M3:M5").Formula = "=L3/1000"

# Match column M's width to column L/N's (31.5703125 chars); nearest value the
# pixel-rounded ColumnWidth model lands on.
$ws.Columns.Item(13).ColumnWidth = 30.7

# Update the active selection / view state
[void]$ws.Range("F1").Select()
$excel.ActiveWindow.ScrollColumn = 6
[void]$ws.Range("M13").Select()

# Update the workbook window position
$wb.Windows.Item(1).Left = -4815
$wb.Windows.Item(1).Top = 2115

$wb.Save()
